# Add a new row (row 9) for a 4th "Sud-Ouest lausannois" parish, matching
# the formatting of the row above it (row 8).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry over the formatting of the last data row onto the new one.
$ws.Range("A8:G8").Copy() | Out-Null
$ws.Range("A9:G9").PasteSpecial(-4122) | Out-Null # xlPasteFormats

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 4070000000
$ws.Range("C9").Value = 4070
$ws.Range("D9").Value = 4000
$ws.Range("E9").Value = "Sud-Ouest lausannois"
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = "P"

# Update the active selection to reflect where Excel was left after the edit.
$ws.Range("E9").Select() | Out-Null
